$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '76.556.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.68%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.042.55'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.40%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '201.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '632.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.04%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  +0.77%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.211'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.19%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.041.16'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.39%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.437'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.44%  '

$ws.Range("E12").Value = '  -0.35%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.11'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.598.54'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.34%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.55'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +6.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.504.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.71%  '

$ws.Range("E17").Value = '  +2.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.014.49'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.45%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.16%  '

$ws.Range("E20").Value = '  +3.45%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.75%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.85%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.09%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.184.03'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.39'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.70%  '

$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.98'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000113'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.32'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.89%  '

$ws.Range("E32").Value = '  +1.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '517.77'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.76%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.98'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.09%  '

$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("E36").Value = '  +3.49%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.58'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.67%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.385'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.70%  '

$ws.Range("E39").Value = '  +1.98%  '

$ws.Range("E40").Value = '  +2.94%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '187.87'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.112'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.71%  '

$ws.Range("E43").Value = '  +0.25%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.17'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.75%  '

$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.26'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +6.02%  '

$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.68'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.47'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.16%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.732'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +11.43%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.609'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.90'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.78%  '
